# Adds a new "ODI Bowling Extra" worksheet (scraped bowling attributes)
# after the existing "ODI Batting Extra" sheet, and tidies up the blank
# placeholder cells left behind on "ODI Batting Extra" by the scraper.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Clean up now-unused blank placeholder cells on "ODI Batting Extra".
#    These rows only ever had a MATCH_CODE (col A) and MAN_OF_MATCH
#    (col F) populated - the B:E scrape columns were written out blank.
# ---------------------------------------------------------------------
$wsBattingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$blankRows = @(3, 5, 6, 8, 9, 11, 18)
foreach ($r in $blankRows) {
    $wsBattingExtra.Range("B" + $r + ":E" + $r).ClearContents()
}
$wsBattingExtra.Range("E21").ClearContents()

# ---------------------------------------------------------------------
# 2) Add the new "ODI Bowling Extra" sheet as the last tab.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# Match the header formatting used on the other "Extra" sheet (bold,
# bordered, centered) by copying its format over.
$wsBattingExtra.Range("A1:C1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "MAIDEN_OVERS"
$ws.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Write MATCH_CODE (col A) as text so numeric-looking codes don't get
# silently converted to numbers.
$ws.Range("A2:C21").NumberFormat = "@"

$data = @(
    @(2,  "3532", $null, $null),
    @(3,  "3533", "0",   "20.00%"),
    @(4,  "3535", $null, $null),
    @(5,  "3538", "0",   $null),
    @(6,  "3580", $null, $null),
    @(7,  "3581", $null, $null),
    @(8,  "3583", "0",   "20.00%"),
    @(9,  "3593", $null, $null),
    @(10, "3596", $null, $null),
    @(11, "3597", "0",   "20.00%"),
    @(12, "3598", $null, $null),
    @(13, "3617", "0",   "30.00%"),
    @(14, "3622", "0",   "10.00%"),
    @(15, "3625", "1",   "20.00%"),
    @(16, "3629", "0",   "30.00%"),
    @(17, "3655", "0",   "40.00%"),
    @(18, "3661", $null, $null),
    @(19, "3678", "0",   "20.00%"),
    @(20, "3680", "0",   "10.00%"),
    @(21, "3683", "0",   $null)
)

foreach ($row in $data) {
    $r = $row[0]
    $matchCode = $row[1]
    $maidenOvers = $row[2]
    $pctWickets = $row[3]

    $ws.Range("A" + $r).Value = $matchCode
    if ($maidenOvers -ne $null) {
        $ws.Range("B" + $r).Value = $maidenOvers
    }
    if ($pctWickets -ne $null) {
        $ws.Range("C" + $r).Value = $pctWickets
    }
}
